# Auto-generated script to apply cryptos.xlsx diff changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.000.89"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "3.078.06"
$ws.Range("E3").Value = "  -1.07%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.65"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.05"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -2.83%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.073.88"
$ws.Range("E8").Value = "  -1.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.514"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.39"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("E11").Value = "  -1.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.471"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -1.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000241"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -1.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.05"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -2.71%  "
$ws.Range("E15").Value = "  -2.02%  "
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "3.590.53"
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "66.950.03"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("E18").Value = "  -1.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.64"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +1.28%  "
$ws.Range("D20").Value = "3.077.65"
$ws.Range("E20").Value = "  -1.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "491.26"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +3.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.687"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -3.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.69"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -2.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.78"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -1.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.86"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -4.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.21"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -2.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.21"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +3.21%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.86"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -0.68%  "
$ws.Range("E30").Value = "  -5.47%  "
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.81"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -2.87%  "
$ws.Range("E33").Value = "  -2.55%  "
$ws.Range("D34").Value = "0.0₃0905"
$ws.Range("E34").Value = "  -2.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.68"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -2.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.954"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -2.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "46.66"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("E39").Value = "  +0.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.98"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -4.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.302"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -2.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.31"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -3.21%  "
$ws.Range("D43").Value = "2.774.82"
$ws.Range("E43").Value = "  -0.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "368.81"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -2.46%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "135.82"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0344"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -2.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.45"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -3.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.59"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -1.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.15"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -2.07%  "
$ws.Range("E51").Value = "  -1.47%  "
